# CompareQuiresList.xlsx - "Pushed - 2 on 20/Oct/2019" edit
#
# 1. Rename the backup sheet "Sheet3" -> "Bkp"
# 2. Update the Master sheet: PROCESS_FLAG values both become "Y", and the
#    TABLE_NAME values for the CNTRL/RISK rows are swapped
#    (MIG_MS_GRC_RISK <-> MIG_MS_GRC_CONTROL)
# 3. Update each sheet's remembered selection / active tab to match the
#    state the workbook was saved in.

$wb = $excel.ActiveWorkbook

$master = $wb.Worksheets.Item("Master")
$tableList = $wb.Worksheets.Item("TableList")
$bkp = $wb.Worksheets.Item("Sheet3")

# --- rename backup sheet ---
$bkp.Name = "Bkp"

# --- data edits on Master ---
$master.Range("D2").Value = "Y"
$master.Range("C3").Value = "MIG_MS_GRC_CONTROL"
$master.Range("C4").Value = "MIG_MS_GRC_RISK"
$master.Range("D5").Value = "Y"

# --- selections / active sheet (TableList ends up active/selected) ---
$master.Range("D7").Select()
$bkp.Range("G16").Select()
$tableList.Range("A6:E9").Select()
